$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 5546.077
$ws.Range("I113").Value = 3419.8
$ws.Range("J113").Value = 6875
$ws.Range("K113").Value = 3419.8
$ws.Range("L113").Value = 6875
$ws.Range("M113").Value = -165.8000000000002
$ws.Range("N113").Value = -13383
# Row 115
$ws.Range("H115").Value = 1240.4166
$ws.Range("I115").Value = 1281.3636
$ws.Range("J115").Value = 790
$ws.Range("K115").Value = 3844.0908
$ws.Range("L115").Value = 2370
$ws.Range("M115").Value = -2277.0908
$ws.Range("N115").Value = -5504
# Row 137
$ws.Range("H137").Value = 1833834.4
$ws.Range("I137").Value = 2382134.8
$ws.Range("J137").Value = 6166.6665
$ws.Range("K137").Value = 7146404.399999999
$ws.Range("L137").Value = 18499.9995
$ws.Range("M137").Value = -7143854.399999999
$ws.Range("N137").Value = -23599.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 876.25
$ws.Range("I2").Value = 876.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 876.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -763.25
$ws.Range("N2").ClearContents()
# Row 61
$ws.Range("H61").Value = 2147.7896
$ws.Range("I61").Value = 1348.75
$ws.Range("K61").Value = 1348.75
$ws.Range("M61").Value = -1136.75
# Row 74
$ws.Range("H74").Value = 5984.05
$ws.Range("I74").Value = 7191.3076
$ws.Range("J74").Value = 3742
$ws.Range("K74").Value = 7191.3076
$ws.Range("L74").Value = 3742
$ws.Range("M74").Value = -6317.3076
$ws.Range("N74").Value = -5490
# Row 77
$ws.Range("H77").Value = 5984.05
$ws.Range("I77").Value = 7191.3076
$ws.Range("J77").Value = 3742
$ws.Range("K77").Value = 35956.538
$ws.Range("L77").Value = 18710
$ws.Range("M77").Value = -31588.538
$ws.Range("N77").Value = -27446
# Row 81
$ws.Range("H81").Value = 40600
$ws.Range("I81").Value = 39000
$ws.Range("J81").Value = 41400
$ws.Range("K81").Value = 39000
$ws.Range("L81").Value = 41400
$ws.Range("M81").Value = -38002
$ws.Range("N81").Value = -43396
# Row 84
$ws.Range("H84").Value = 40600
$ws.Range("I84").Value = 39000
$ws.Range("J84").Value = 41400
$ws.Range("K84").Value = 117000
$ws.Range("L84").Value = 124200
$ws.Range("M84").Value = -112008
$ws.Range("N84").Value = -134184
# Row 88
$ws.Range("H88").Value = 7411407
$ws.Range("I88").Value = 16670041
$ws.Range("J88").Value = 4500
$ws.Range("K88").Value = 16670041
$ws.Range("L88").Value = 4500
$ws.Range("M88").Value = -16669635
$ws.Range("N88").Value = -5312
# Row 91
$ws.Range("H91").Value = 7411407
$ws.Range("I91").Value = 16670041
$ws.Range("J91").Value = 4500
$ws.Range("K91").Value = 16670041
$ws.Range("L91").Value = 4500
$ws.Range("M91").Value = -16668637
$ws.Range("N91").Value = -7308
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
# Row 116
$ws.Range("H116").Value = 876.25
$ws.Range("I116").Value = 876.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 876.25
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1417.75
$ws.Range("N116").ClearContents()
# Row 122
$ws.Range("H122").Value = 2642.2
$ws.Range("I122").Value = 1741
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 5223
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -2773
$ws.Range("N122").Value = -30400
# Row 132
$ws.Range("H132").Value = 2649.889
$ws.Range("I132").Value = 1144.5454
$ws.Range("J132").Value = 5015.4287
$ws.Range("K132").Value = 3433.6362
$ws.Range("L132").Value = 15046.2861
$ws.Range("M132").Value = -903.6361999999999
$ws.Range("N132").Value = -20106.2861
# Row 136
$ws.Range("H136").Value = 2147.7896
$ws.Range("I136").Value = 1348.75
$ws.Range("K136").Value = 4046.25
$ws.Range("M136").Value = -1496.25
# Row 137
$ws.Range("H137").Value = 39968.332
$ws.Range("J137").Value = 39968.332
$ws.Range("L137").Value = 39968.332
$ws.Range("N137").Value = -50168.332

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 876.25
$ws.Range("I3").Value = 876.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 876.25
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -762.25
$ws.Range("N3").ClearContents()
# Row 21
$ws.Range("H21").Value = 37542
$ws.Range("J21").Value = 37542
$ws.Range("L21").Value = 37542
$ws.Range("N21").Value = -38014
# Row 70
$ws.Range("H70").Value = 71300
$ws.Range("J70").Value = 71300
$ws.Range("L70").Value = 71300
$ws.Range("N70").Value = -71886
# Row 73
$ws.Range("H73").Value = 71300
$ws.Range("J73").Value = 71300
$ws.Range("L73").Value = 71300
$ws.Range("N73").Value = -73328
# Row 86
$ws.Range("H86").Value = 2150
$ws.Range("I86").Value = 2080
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 2080
$ws.Range("L86").Value = 2200
$ws.Range("M86").Value = -957
$ws.Range("N86").Value = -4446
# Row 89
$ws.Range("H89").Value = 2150
$ws.Range("I89").Value = 2080
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 10400
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = -4784
$ws.Range("N89").Value = -22232
# Row 134
$ws.Range("H134").Value = 1885.3334
$ws.Range("I134").Value = 1282.7333
$ws.Range("J134").Value = 4898.3335
$ws.Range("K134").Value = 3848.199900000001
$ws.Range("L134").Value = 14695.0005
$ws.Range("M134").Value = -1313.199900000001
$ws.Range("N134").Value = -19765.0005
# Row 137
$ws.Range("H137").Value = 41206.668
$ws.Range("J137").Value = 41206.668
$ws.Range("L137").Value = 41206.668
$ws.Range("N137").Value = -51406.668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 25545.715
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 34964
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 34964
$ws.Range("M23").Value = -1760
$ws.Range("N23").Value = -35444
# Row 27
$ws.Range("H27").Value = 25545.715
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 34964
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 34964
$ws.Range("M27").Value = -1808
$ws.Range("N27").Value = -35348
# Row 31
$ws.Range("H31").Value = 2859.838
$ws.Range("I31").Value = 1022.7647
$ws.Range("J31").Value = 4421.35
$ws.Range("K31").Value = 1022.7647
$ws.Range("L31").Value = 4421.35
$ws.Range("M31").Value = -727.7646999999999
$ws.Range("N31").Value = -5011.35
# Row 34
$ws.Range("H34").Value = 2859.838
$ws.Range("I34").Value = 1022.7647
$ws.Range("J34").Value = 4421.35
$ws.Range("K34").Value = 1022.7647
$ws.Range("L34").Value = 4421.35
$ws.Range("M34").Value = -820.7646999999999
$ws.Range("N34").Value = -4825.35
# Row 58
$ws.Range("H58").Value = 2482.92
$ws.Range("I58").Value = 1570.5088
$ws.Range("J58").Value = 5372.222
$ws.Range("K58").Value = 1570.5088
$ws.Range("L58").Value = 5372.222
$ws.Range("M58").Value = -1367.5088
$ws.Range("N58").Value = -5778.222
# Row 132
$ws.Range("H132").Value = 2740.9546
$ws.Range("I132").Value = 1612.6875
$ws.Range("J132").Value = 5749.6665
$ws.Range("K132").Value = 4838.0625
$ws.Range("L132").Value = 17248.9995
$ws.Range("M132").Value = -2308.0625
$ws.Range("N132").Value = -22308.9995
# Row 134
$ws.Range("H134").Value = 9212.134
$ws.Range("I134").Value = 13723
$ws.Range("J134").Value = 4056.8572
$ws.Range("K134").Value = 41169
$ws.Range("L134").Value = 12170.5716
$ws.Range("M134").Value = -38634
$ws.Range("N134").Value = -17240.5716
# Row 136
$ws.Range("H136").Value = 2482.92
$ws.Range("I136").Value = 1570.5088
$ws.Range("J136").Value = 5372.222
$ws.Range("K136").Value = 4711.526400000001
$ws.Range("L136").Value = 16116.666
$ws.Range("M136").Value = -2161.526400000001
$ws.Range("N136").Value = -21216.666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 3572011.8
$ws.Range("I113").Value = 603.8889
$ws.Range("J113").Value = 7353502.5
$ws.Range("K113").Value = 1811.6667
$ws.Range("L113").Value = 22060507.5
$ws.Range("M113").Value = 358.3332999999998
$ws.Range("N113").Value = -22064847.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4325.5625
$ws.Range("I132").Value = 2338.875
$ws.Range("J132").Value = 6312.25
$ws.Range("K132").Value = 7016.625
$ws.Range("L132").Value = 18936.75
$ws.Range("M132").Value = -4486.625
$ws.Range("N132").Value = -23996.75
# Row 137
$ws.Range("H137").Value = 37185
$ws.Range("J137").Value = 37185
$ws.Range("L137").Value = 37185
$ws.Range("N137").Value = -47385

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 3250

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 93
$ws.Range("H93").Value = 39800
$ws.Range("J93").Value = 39800
$ws.Range("L93").Value = 39800
$ws.Range("N93").Value = -44792
# Row 122
$ws.Range("H122").Value = 5909.273
$ws.Range("J122").Value = 7714.2856
$ws.Range("L122").Value = 23142.8568
$ws.Range("N122").Value = -28042.8568
# Row 136
$ws.Range("H136").Value = 4670.174
$ws.Range("I136").Value = 2322.2144
$ws.Range("J136").Value = 8322.556
$ws.Range("K136").Value = 6966.6432
$ws.Range("L136").Value = 24967.668
$ws.Range("M136").Value = -4416.6432
$ws.Range("N136").Value = -30067.668

Write-Output "Applied all updates."